$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("core i9 13900k")

# Add new data row for version 7.1.0
$ws.Range("A9").Value = "7.1.0"
$ws.Range("B9").Value = 27.75
$ws.Range("C9").Value = 25746866156
$ws.Range("D9").Value = 927814996

# Extend the three line-chart series (Time, Nodes, NPS) to include the new row
$co1 = $ws.ChartObjects().Item(1)
$s1 = $co1.Chart.SeriesCollection().Item(1)
$s1.Formula = "=SERIES('core i9 13900k'!`$B`$1,'core i9 13900k'!`$A`$2:`$A`$36,'core i9 13900k'!`$B`$2:`$B`$9,1)"

$co2 = $ws.ChartObjects().Item(2)
$s2 = $co2.Chart.SeriesCollection().Item(1)
$s2.Formula = "=SERIES('core i9 13900k'!`$C`$1,'core i9 13900k'!`$A`$2:`$A`$36,'core i9 13900k'!`$C`$2:`$C`$9,1)"

$co3 = $ws.ChartObjects().Item(3)
$s3 = $co3.Chart.SeriesCollection().Item(1)
$s3.Formula = "=SERIES('core i9 13900k'!`$D`$1,'core i9 13900k'!`$A`$2:`$A`$36,'core i9 13900k'!`$D`$2:`$D`$9,1)"

# Move the active selection as in the source workbook
$ws.Range("D17").Select()
